$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 ("PWM", 5, "Shooter Angle") was a leftover/duplicate data row - remove it entirely.
$ws.Rows.Item(18).Delete()

# After the delete above, the former rows 24/25 (the separate Left/Right
# "pulling piston in" Super Shifter solenoid rows) are now rows 23/24.
# Collapse them into a single generic "Super Shifter" row and drop the duplicate.
$ws.Range("A23").Value = "DRIVE_SUPERSHIFTSOLIN"
$ws.Range("D23").Value = "Super Shifter solenoid for pulling piston in"
$ws.Rows.Item(24).Delete()

# The former rows 26/27 (the separate Left/Right "pushing piston out" Super
# Shifter solenoid rows) are now row 24/25. Collapse them as well.
$ws.Range("A24").Value = "DRIVE_SUPERSHIFTSOLOUT"
$ws.Range("D24").Value = "Super Shifter solenoid for pushing piston out"
$ws.Rows.Item(25).Delete()
